$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-11 17:18:37"
$ws.Range("H2").Value = "'78%"
$ws.Range("K2").Value = "5.3 MJ/m2"
$ws.Range("E3").Value = "2026-02-11 17:18:39"
$ws.Range("O3").Value = "0.0 °C"
$ws.Range("E4").Value = "2026-02-11 17:18:42"
$ws.Range("H4").Value = "'57%"
$ws.Range("J4").Value = "1002.9 hPa"
$ws.Range("E5").Value = "2026-02-11 17:18:44"
$ws.Range("O5").Value = "0.3 °C"
$ws.Range("E6").Value = "2026-02-11 17:18:47"
$ws.Range("J6").Value = "1003.4 hPa"
$ws.Range("O6").Value = "13.2 °C"
$ws.Range("E7").Value = "2026-02-11 17:18:50"
$ws.Range("J7").Value = "1003.8 hPa"
$ws.Range("E8").Value = "2026-02-11 17:18:52"
$ws.Range("H8").Value = "'52%"
$ws.Range("E9").Value = "2026-02-11 17:18:54"
$ws.Range("E10").Value = "2026-02-11 17:18:57"
$ws.Range("E11").Value = "2026-02-11 17:18:59"
$ws.Range("H11").Value = "'80%"
$ws.Range("O11").Value = "7.8 °C"
$ws.Range("E12").Value = "2026-02-11 17:19:02"
$ws.Range("E13").Value = "2026-02-11 17:19:04"
$ws.Range("J13").Value = "1005.5 hPa"
$ws.Range("O13").Value = "7.6 °C"
$ws.Range("E14").Value = "2026-02-11 17:19:07"
$ws.Range("K14").Value = "10.4 MJ/m2"
$ws.Range("E15").Value = "2026-02-11 17:19:10"
$ws.Range("E16").Value = "2026-02-11 17:19:12"
$ws.Range("H16").Value = "'62%"
$ws.Range("I16").Value = "4.6 mm"
$ws.Range("K16").Value = "7.4 MJ/m2"
$ws.Range("E17").Value = "2026-02-11 17:19:14"
$ws.Range("E18").Value = "2026-02-11 17:19:17"
$ws.Range("J18").Value = "1003.4 hPa"
$ws.Range("K18").Value = "10.8 MJ/m2"
$ws.Range("O18").Value = "13.8 °C"
$ws.Range("E19").Value = "2026-02-11 17:19:20"
$ws.Range("K19").Value = "10.3 MJ/m2"
$ws.Range("E20").Value = "2026-02-11 17:19:22"
$ws.Range("E21").Value = "2026-02-11 17:19:25"
$ws.Range("J21").Value = "1006.0 hPa"
$ws.Range("K21").Value = "9.2 MJ/m2"
$ws.Range("E22").Value = "2026-02-11 17:19:27"
$ws.Range("I22").Value = "1.8 mm"
$ws.Range("O22").Value = "-2.8 °C"
$ws.Range("E23").Value = "2026-02-11 17:19:30"
$ws.Range("I23").Value = "3.8 mm"
$ws.Range("E24").Value = "2026-02-11 17:19:32"
$ws.Range("H24").Value = "'71%"
$ws.Range("I24").Value = "5.4 mm"
$ws.Range("J24").Value = "1007.5 hPa"
$ws.Range("N24").Value = "11.3 °C 16:59 TU"
$ws.Range("E25").Value = "2026-02-11 17:19:35"
$ws.Range("H25").Value = "'61%"
$ws.Range("L25").Value = "44.3 km/h - 244º 16:57 TU"
$ws.Range("E26").Value = "2026-02-11 17:19:37"
$ws.Range("J26").Value = "1003.3 hPa"
$ws.Range("K26").Value = "10.7 MJ/m2"
$ws.Range("O26").Value = "7.0 °C"
$ws.Range("E27").Value = "2026-02-11 17:19:40"
$ws.Range("H27").Value = "'81%"
$ws.Range("E28").Value = "2026-02-11 17:19:43"
$ws.Range("H28").Value = "'82%"
$ws.Range("J28").Value = "1003.7 hPa"
$ws.Range("O28").Value = "10.7 °C"
$ws.Range("E29").Value = "2026-02-11 17:19:45"
$ws.Range("K29").Value = "10.6 MJ/m2"
$ws.Range("O29").Value = "12.8 °C"
$ws.Range("E30").Value = "2026-02-11 17:19:48"
$ws.Range("J30").Value = "1003.6 hPa"
$ws.Range("K30").Value = "10.6 MJ/m2"
$ws.Range("E31").Value = "2026-02-11 17:19:50"
$ws.Range("J31").Value = "1002.8 hPa"
$ws.Range("K31").Value = "9.5 MJ/m2"
$ws.Range("E32").Value = "2026-02-11 17:19:53"
$ws.Range("H32").Value = "'73%"
$ws.Range("I32").Value = "2.9 mm"
$ws.Range("E33").Value = "2026-02-11 17:19:56"
$ws.Range("H33").Value = "'80%"
$ws.Range("K33").Value = "10.9 MJ/m2"
$ws.Range("E34").Value = "2026-02-11 17:19:58"
$ws.Range("E35").Value = "2026-02-11 17:20:01"
$ws.Range("G35").Value = "1 cm"
$ws.Range("H35").Value = "'68%"
$ws.Range("J35").Value = "1008.3 hPa"
$ws.Range("E36").Value = "2026-02-11 17:20:03"
$ws.Range("J36").Value = "1003.8 hPa"
$ws.Range("K36").Value = "10.6 MJ/m2"
$ws.Range("O36").Value = "12.9 °C"
$ws.Range("E37").Value = "2026-02-11 17:20:06"
$ws.Range("J37").Value = "1004.9 hPa"
$ws.Range("O37").Value = "9.3 °C"
$ws.Range("E38").Value = "2026-02-11 17:20:08"
$ws.Range("H38").Value = "'59%"
$ws.Range("O38").Value = "15.7 °C"
$ws.Range("E39").Value = "2026-02-11 17:20:11"
$ws.Range("O39").Value = "1.0 °C"
$ws.Range("E40").Value = "2026-02-11 17:20:14"
$ws.Range("I40").Value = "0.7 mm"
$ws.Range("J40").Value = "1007.4 hPa"
$ws.Range("O40").Value = "7.5 °C"
$ws.Range("E41").Value = "2026-02-11 17:20:16"
$ws.Range("H41").Value = "'43%"
$ws.Range("K41").Value = "9.0 MJ/m2"
$ws.Range("N41").Value = "16.9 °C 16:48 TU"
$ws.Range("O41").Value = "19.4 °C"
$ws.Range("E42").Value = "2026-02-11 17:20:19"
$ws.Range("E43").Value = "2026-02-11 17:20:21"
$ws.Range("E44").Value = "2026-02-11 17:20:24"
$ws.Range("H44").Value = "'82%"
$ws.Range("I44").Value = "4.3 mm"
$ws.Range("E45").Value = "2026-02-11 17:20:27"
$ws.Range("H45").Value = "'86%"
$ws.Range("J45").Value = "1006.4 hPa"
$ws.Range("E46").Value = "2026-02-11 17:20:29"
$ws.Range("H46").Value = "'54%"
$ws.Range("J46").Value = "1007.9 hPa"
$ws.Range("K46").Value = "7.5 MJ/m2"
$ws.Range("N46").Value = "13.7 °C 16:59 TU"
$ws.Range("O46").Value = "17.6 °C"
